{"js": "/*\n * Change 1: Split the \"M\u00e9todo\" run's text at the sentence boundary between\n *   \"...exerc\u00edcios dirigidos. \" and \"Avalia\u00e7\u00e3o baseada...\" by inserting a\n *   manual line break (rendered as <w:br/> in OOXML) at that point.\n *\n * Change 2: Split the big Bibliografia run into one run-of-text per\n *   reference, separated by manual line breaks, with THREE line breaks\n *   before the \"Bibliografia complementar:\" heading.\n *\n * Both edits are performed with Body.search() to locate the exact text and\n * Range.insertText(..., \"Replace\") to rewrite it in place; embedding a\n * vertical-tab character (\\v / U+000B) in the replacement text is how\n * Office.js represents an explicit line break, which the host serializes\n * back out as <w:br/>.\n */\n\nconst body = context.document.body;\n\n// --- Change 1: \"M\u00e9todo\" sentence -> add a line break before \"Avalia\u00e7\u00e3o baseada\" ---\nconst original1 = \"Aulas te\u00f3ricas e pr\u00e1ticas, visitas t\u00e9cnicas e exerc\u00edcios dirigidos. Avalia\u00e7\u00e3o baseada em provas, exerc\u00edcios e trabalhos pr\u00e1ticos e relat\u00f3rios.\";\nconst replacement1 = \"Aulas te\u00f3ricas e pr\u00e1ticas, visitas t\u00e9cnicas e exerc\u00edcios dirigidos. \\u000bAvalia\u00e7\u00e3o baseada em provas, exerc\u00edcios e trabalhos pr\u00e1ticos e relat\u00f3rios.\";\n\nconst results1 = body.search(original1, { matchCase: true });\nresults1.load(\"items\");\nawait context.sync();\n\nif (results1.items.length === 0) {\n  throw new Error(\"Change 1: target text not found\");\n}\nresults1.items[0].insertText(replacement1, \"Replace\");\nawait context.sync();\n\n// --- Change 2: Bibliografia paragraph -> one reference per line ---\nconst original2 = \"BURROUGH, P. A. Principles of Geographical Information Systems - Spatial Information Systems and Geoestatistics, Oxford: Clarendon Press, 1998.BURROUGH, P. A.; MCDONNELL, R. A. Principles of Geographical Information Systems. Oxford University Press, 1998.C\u00c2MARA, G. & MEDEIROS, J. S. GIS para Meio Ambiente. INPE. S\u00e3o Jos\u00e9 dos Campos, SP, 1998.CROSTA, A. P. Processamento Digital de Imagens de Sensoriamento Remoto. Campinas \u2013 SP, 1992.FLORENZANO, T. G. Imagens de Sat\u00e9lite para Estudos Ambientais. Oficina de textos. S\u00e3o Paulo, 2002.IBGE. No\u00e7\u00f5es B\u00e1sicas de Cartografia. Rio de Janeiro. Cole\u00e7\u00e3o Manuais T\u00e9cnicos em Geoci\u00eancias, 1999.LONGLEY, P. A.; GOODCHILD, M. F.; MAGUIRE, D. J.; RHIND, D. W. Geographic Information Systems and Science. John Wiley & Sons, 2001.MIRANDA, J. I.; Fundamentos de Sistemas de Informa\u00e7\u00f5es Geogr\u00e1ficas. Bras\u00edlia, Embrapa, 2005.MOREIRA, M. A. Fundamentos do Sensoriamento Remoto e Metodologias de Aplica\u00e7\u00e3o. S\u00e3o Jos\u00e9 dos Campos \u2013 SP \u2013 INPE, 2001.SILVA, A.B. Sistemas de Informa\u00e7\u00f5es Geo-referenciadas. Editora da Unicamp. Campinas. 1999.SILVA, A. B; Sistemas de informa\u00e7\u00f5es Geo-referenciadas: conceitos e fundamentos. Campinas: Editora da Unicamp, 2003.SILVA, J.X. Geoprocessamento para An\u00e1lise Ambiental. Rio de Janeiro. 2001.Bibliografia complementar:CARVALHO, M. S.; PINA, M. F.; SANTOS, S. M.  Conceitos B\u00e1sicos de Sistemas de Informa\u00e7\u00e3o Geogr\u00e1fica e Cartografia Aplicados \u00e0 Sa\u00fade. Rede Interagencial de Informa\u00e7\u00f5es para a Sa\u00fade. Bras\u00edlia. Minist\u00e9rio da Sa\u00fade, 2000.DENT, B. D.  Cartography Thematic Map Design. 5th Edition. WCB/McGraw-Hill, 1999.MATOS, J. Fundamentos da Informa\u00e7\u00e3o Geogr\u00e1fica. Lisboa, Lidel, 2008.MORAES NOVO, E. M. L. Sensoriamento Remoto \u2013 Princ\u00edpios e Aplica\u00e7\u00f5es. 2\u00aaEdi\u00e7\u00e3o. S\u00e3o Paulo, 1992.\";\nconst replacement2 = \"BURROUGH, P. A. Principles of Geographical Information Systems - Spatial Information Systems and Geoestatistics, Oxford: Clarendon Press, 1998.\\u000bBURROUGH, P. A.; MCDONNELL, R. A. Principles of Geographical Information Systems. Oxford University Press, 1998.\\u000bC\u00c2MARA, G. & MEDEIROS, J. S. GIS para Meio Ambiente. INPE. S\u00e3o Jos\u00e9 dos Campos, SP, 1998.\\u000bCROSTA, A. P. Processamento Digital de Imagens de Sensoriamento Remoto. Campinas \u2013 SP, 1992.\\u000bFLORENZANO, T. G. Imagens de Sat\u00e9lite para Estudos Ambientais. Oficina de textos. S\u00e3o Paulo, 2002.\\u000bIBGE. No\u00e7\u00f5es B\u00e1sicas de Cartografia. Rio de Janeiro. Cole\u00e7\u00e3o Manuais T\u00e9cnicos em Geoci\u00eancias, 1999.\\u000bLONGLEY, P. A.; GOODCHILD, M. F.; MAGUIRE, D. J.; RHIND, D. W. Geographic Information Systems and Science. John Wiley & Sons, 2001.\\u000bMIRANDA, J. I.; Fundamentos de Sistemas de Informa\u00e7\u00f5es Geogr\u00e1ficas. Bras\u00edlia, Embrapa, 2005.\\u000bMOREIRA, M. A. Fundamentos do Sensoriamento Remoto e Metodologias de Aplica\u00e7\u00e3o. S\u00e3o Jos\u00e9 dos Campos \u2013 SP \u2013 INPE, 2001.\\u000bSILVA, A.B. Sistemas de Informa\u00e7\u00f5es Geo-referenciadas. Editora da Unicamp. Campinas. 1999.\\u000bSILVA, A. B; Sistemas de informa\u00e7\u00f5es Geo-referenciadas: conceitos e fundamentos. Campinas: Editora da Unicamp, 2003.\\u000bSILVA, J.X. Geoprocessamento para An\u00e1lise Ambiental. Rio de Janeiro. 2001.\\u000b\\u000b\\u000bBibliografia complementar:\\u000bCARVALHO, M. S.; PINA, M. F.; SANTOS, S. M.  Conceitos B\u00e1sicos de Sistemas de Informa\u00e7\u00e3o Geogr\u00e1fica e Cartografia Aplicados \u00e0 Sa\u00fade. Rede Interagencial de Informa\u00e7\u00f5es para a Sa\u00fade. Bras\u00edlia. Minist\u00e9rio da Sa\u00fade, 2000.\\u000bDENT, B. D.  Cartography Thematic Map Design. 5th Edition. WCB/McGraw-Hill, 1999.\\u000bMATOS, J. Fundamentos da Informa\u00e7\u00e3o Geogr\u00e1fica. Lisboa, Lidel, 2008.\\u000bMORAES NOVO, E. M. L. Sensoriamento Remoto \u2013 Princ\u00edpios e Aplica\u00e7\u00f5es. 2\u00aaEdi\u00e7\u00e3o. S\u00e3o Paulo, 1992.\";\n\nconst results2 = body.search(original2, { matchCase: true });\nresults2.load(\"items\");\nawait context.sync();\n\nif (results2.items.length === 0) {\n  throw new Error(\"Change 2: target text not found\");\n}\nresults2.items[0].insertText(replacement2, \"Replace\");\nawait context.sync();\n", "ps1": "# Applies the LOB1214.docx edit described in the commit diff:\n#   1. In the \"M\u00e9todo\" sentence (Avalia\u00e7\u00e3o section), insert a manual line\n#      break between \"...exerc\u00edcios dirigidos. \" and \"Avalia\u00e7\u00e3o baseada...\".\n#   2. In the Bibliografia paragraph, put each reference on its own line\n#      (manual line breaks), with THREE line breaks before the\n#      \"Bibliografia complementar:\" heading.\n#\n# Both edits use Range.Find/Replacement with wildcards enabled so the\n# \"^l\" token in the replacement text is expanded to a manual line break\n# (OOXML <w:br/>), matching the target diff exactly.\n\n$d = $word.ActiveDocument\n\n# --- Change 1 --------------------------------------------------------\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = 'Aulas te\u00f3ricas e pr\u00e1ticas, visitas t\u00e9cnicas e exerc\u00edcios dirigidos. Avalia\u00e7\u00e3o baseada em provas, exerc\u00edcios e trabalhos pr\u00e1ticos e relat\u00f3rios.'\n$find1.Replacement.Text = 'Aulas te\u00f3ricas e pr\u00e1ticas, visitas t\u00e9cnicas e exerc\u00edcios dirigidos. ^lAvalia\u00e7\u00e3o baseada em provas, exerc\u00edcios e trabalhos pr\u00e1ticos e relat\u00f3rios.'\n$ok1 = $find1.Execute(\n    $find1.Text, $false, $false, $false, $false, $false, $true, 1, $false,\n    $find1.Replacement.Text, 2\n)\nif (-not $ok1) {\n    throw \"Change 1: search text not found\"\n}\n\n# --- Change 2 --------------------------------------------------------\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = 'BURROUGH, P. A. Principles of Geographical Information Systems - Spatial Information Systems and Geoestatistics, Oxford: Clarendon Press, 1998.BURROUGH, P. A.; MCDONNELL, R. A. Principles of Geographical Information Systems. Oxford University Press, 1998.C\u00c2MARA, G. & MEDEIROS, J. S. GIS para Meio Ambiente. INPE. S\u00e3o Jos\u00e9 dos Campos, SP, 1998.CROSTA, A. P. Processamento Digital de Imagens de Sensoriamento Remoto. Campinas \u2013 SP, 1992.FLORENZANO, T. G. Imagens de Sat\u00e9lite para Estudos Ambientais. Oficina de textos. S\u00e3o Paulo, 2002.IBGE. No\u00e7\u00f5es B\u00e1sicas de Cartografia. Rio de Janeiro. Cole\u00e7\u00e3o Manuais T\u00e9cnicos em Geoci\u00eancias, 1999.LONGLEY, P. A.; GOODCHILD, M. F.; MAGUIRE, D. J.; RHIND, D. W. Geographic Information Systems and Science. John Wiley & Sons, 2001.MIRANDA, J. I.; Fundamentos de Sistemas de Informa\u00e7\u00f5es Geogr\u00e1ficas. Bras\u00edlia, Embrapa, 2005.MOREIRA, M. A. Fundamentos do Sensoriamento Remoto e Metodologias de Aplica\u00e7\u00e3o. S\u00e3o Jos\u00e9 dos Campos \u2013 SP \u2013 INPE, 2001.SILVA, A.B. Sistemas de Informa\u00e7\u00f5es Geo-referenciadas. Editora da Unicamp. Campinas. 1999.SILVA, A. B; Sistemas de informa\u00e7\u00f5es Geo-referenciadas: conceitos e fundamentos. Campinas: Editora da Unicamp, 2003.SILVA, J.X. Geoprocessamento para An\u00e1lise Ambiental. Rio de Janeiro. 2001.Bibliografia complementar:CARVALHO, M. S.; PINA, M. F.; SANTOS, S. M.  Conceitos B\u00e1sicos de Sistemas de Informa\u00e7\u00e3o Geogr\u00e1fica e Cartografia Aplicados \u00e0 Sa\u00fade. Rede Interagencial de Informa\u00e7\u00f5es para a Sa\u00fade. Bras\u00edlia. Minist\u00e9rio da Sa\u00fade, 2000.DENT, B. D.  Cartography Thematic Map Design. 5th Edition. WCB/McGraw-Hill, 1999.MATOS, J. Fundamentos da Informa\u00e7\u00e3o Geogr\u00e1fica. Lisboa, Lidel, 2008.MORAES NOVO, E. M. L. Sensoriamento Remoto \u2013 Princ\u00edpios e Aplica\u00e7\u00f5es. 2\u00aaEdi\u00e7\u00e3o. S\u00e3o Paulo, 1992.'\n$find2.Replacement.Text = 'BURROUGH, P. A. Principles of Geographical Information Systems - Spatial Information Systems and Geoestatistics, Oxford: Clarendon Press, 1998.^lBURROUGH, P. A.; MCDONNELL, R. A. Principles of Geographical Information Systems. Oxford University Press, 1998.^lC\u00c2MARA, G. & MEDEIROS, J. S. GIS para Meio Ambiente. INPE. S\u00e3o Jos\u00e9 dos Campos, SP, 1998.^lCROSTA, A. P. Processamento Digital de Imagens de Sensoriamento Remoto. Campinas \u2013 SP, 1992.^lFLORENZANO, T. G. Imagens de Sat\u00e9lite para Estudos Ambientais. Oficina de textos. S\u00e3o Paulo, 2002.^lIBGE. No\u00e7\u00f5es B\u00e1sicas de Cartografia. Rio de Janeiro. Cole\u00e7\u00e3o Manuais T\u00e9cnicos em Geoci\u00eancias, 1999.^lLONGLEY, P. A.; GOODCHILD, M. F.; MAGUIRE, D. J.; RHIND, D. W. Geographic Information Systems and Science. John Wiley & Sons, 2001.^lMIRANDA, J. I.; Fundamentos de Sistemas de Informa\u00e7\u00f5es Geogr\u00e1ficas. Bras\u00edlia, Embrapa, 2005.^lMOREIRA, M. A. Fundamentos do Sensoriamento Remoto e Metodologias de Aplica\u00e7\u00e3o. S\u00e3o Jos\u00e9 dos Campos \u2013 SP \u2013 INPE, 2001.^lSILVA, A.B. Sistemas de Informa\u00e7\u00f5es Geo-referenciadas. Editora da Unicamp. Campinas. 1999.^lSILVA, A. B; Sistemas de informa\u00e7\u00f5es Geo-referenciadas: conceitos e fundamentos. Campinas: Editora da Unicamp, 2003.^lSILVA, J.X. Geoprocessamento para An\u00e1lise Ambiental. Rio de Janeiro. 2001.^l^l^lBibliografia complementar:^lCARVALHO, M. S.; PINA, M. F.; SANTOS, S. M.  Conceitos B\u00e1sicos de Sistemas de Informa\u00e7\u00e3o Geogr\u00e1fica e Cartografia Aplicados \u00e0 Sa\u00fade. Rede Interagencial de Informa\u00e7\u00f5es para a Sa\u00fade. Bras\u00edlia. Minist\u00e9rio da Sa\u00fade, 2000.^lDENT, B. D.  Cartography Thematic Map Design. 5th Edition. WCB/McGraw-Hill, 1999.^lMATOS, J. Fundamentos da Informa\u00e7\u00e3o Geogr\u00e1fica. Lisboa, Lidel, 2008.^lMORAES NOVO, E. M. L. Sensoriamento Remoto \u2013 Princ\u00edpios e Aplica\u00e7\u00f5es. 2\u00aaEdi\u00e7\u00e3o. S\u00e3o Paulo, 1992.'\n$ok2 = $find2.Execute(\n    $find2.Text, $false, $false, $false, $false, $false, $true, 1, $false,\n    $find2.Replacement.Text, 2\n)\nif (-not $ok2) {\n    throw \"Change 2: search text not found\"\n}\n"}
